$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats from the (now-shifted) E column into the new blank D column
# so the new column inherits the same date/number formatting as its neighbours.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting-period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 112900
$ws.Range("D9").Value = 17400
$ws.Range("D10").Value = 95500
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 80200
$ws.Range("D18").Value = 32700
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 35800
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 32700
$ws.Range("D24").Value = 7700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 25000
$ws.Range("D27").Value = 24500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 24500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 24500
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 97200
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 5200
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = "NA"
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 1021400
$ws.Range("D48").Value = 4300
$ws.Range("D49").Value = 15300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 14000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1167000
$ws.Range("D57").Value = 36400
$ws.Range("D58").Value = 829600
$ws.Range("D59").Value = 3800
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 76200
$ws.Range("D62").Value = 22600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 968500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 114900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 198500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 24500
$ws.Range("D83").Value = 3100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 84400
$ws.Range("D91").Value = -724600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -126900
$ws.Range("D96").Value = -6900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 86600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 44100

# Two historical figures were corrected alongside the shift
$ws.Range("F94").Value = -138500
$ws.Range("F102").Value = 1400
